# UiComponentClassDiagram.pptx - "update UI class diagram"
#
# 1) Refresh the cached "datetimeFigureOut" field text (Insert > Header &
#    Footer > Date: "Update automatically") from 10/22/2018 to 10/31/18
#    across every slide layout, the slide master, and the notes master.
# 2) Rename the "BrowserPanel" rectangle's label to "SidebarPanel" on
#    slide 1.

$p = $ppt.ActivePresentation

$newDate = "10/31/18"

# --- 1) Date placeholders -------------------------------------------------

# Slide master: "Date Placeholder 3" is shape index 3.
$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# Notes master: "Date Placeholder 2" is shape index 2.
$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# Every slide layout off the master - the date placeholder's shape index
# varies per layout, so map layout -> shape index explicitly.
$layoutDateShapeIndex = @{
    1 = 3   # Title Slide
    2 = 3   # Title and Content
    3 = 3   # Section Header
    4 = 4   # Two Content
    5 = 6   # Comparison
    6 = 2   # Title Only
    7 = 1   # Blank
    8 = 4   # Content with Caption
    9 = 4   # Picture with Caption
    10 = 3  # Title and Vertical Text
    11 = 3  # Vertical Title and Text
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    $idx = $layoutDateShapeIndex[$i]
    $layout.Shapes.Item($idx).TextFrame.TextRange.Text = $newDate
}

# --- 2) Rename BrowserPanel -> SidebarPanel on slide 1 --------------------

$slide = $p.Slides.Item(1)
$slide.Shapes.Item(9).TextFrame.TextRange.Text = "SidebarPanel"
